$wb = $excel.ActiveWorkbook

# --- Sheet "Generator Data": remove upgrade 2 / upgrade 3 rows (Battery Replacement removed) ---
$ws1 = $wb.Worksheets.Item("Generator Data")

# Delete from bottom to top so row numbers of earlier deletions stay valid
$ws1.Rows("13:14").Delete()   # Yearly O&M Cost at upgrade 2 / 3
$ws1.Rows("10:11").Delete()   # Investment at upgrade 2 / 3
$ws1.Rows("7:8").Delete()     # Nominal Capacity at upgrade 2 / 3

# Update the remaining values that shifted with the recalculated model
$ws1.Cells.Item(6,2).Value = 104752.7053983763     # Nominal Capacity at upgrade 1
$ws1.Cells.Item(7,2).Value = 44006.611537857898    # Investment at upgrade 1
$ws1.Cells.Item(8,2).Value = 4400.6611537857898    # Yearly O&M Cost at upgrade 1
$ws1.Cells.Item(9,2).Value = 737597.48794515361    # Total actualized Fuel Cost

# --- Sheet "Yearly Fuel Costs": add year 4 and year 5 rows, update values ---
$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")

$ws2.Cells.Item(2,2).Value = 147515.0489019602     # Total Fuel Cost at y = 1
$ws2.Cells.Item(3,2).Value = 147520.60976081199    # Total Fuel Cost at y = 2
$ws2.Cells.Item(4,2).Value = 147520.60976081199    # Total Fuel Cost at y = 3

$ws2.Cells.Item(5,1).Value = "Total Fuel Cost at y = 4"
$ws2.Cells.Item(5,2).Value = 147520.60976081199
$ws2.Cells.Item(6,1).Value = "Total Fuel Cost at y = 5"
$ws2.Cells.Item(6,2).Value = 147520.60976081199

# Copy the formatting (style) of row 4's label cell onto the two new label cells
$ws2.Cells.Item(4,1).Copy()
$ws2.Cells.Item(5,1).PasteSpecial(-4122)
$ws2.Cells.Item(6,1).PasteSpecial(-4122)
